$wb = $excel.ActiveWorkbook

# 1) Update the "Status" value from "Ready for handoff" to "In Translation".
#    The same shared text shows up in three places:
#      - Overview!E2 (zh-cn status) and Overview!F2 (de-de status)
#      - zh-cn!C2 (Status column)
#      - de-de!C2 (Status column)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# 2) Narrow a few columns that previously shared the same (wider) width.
#    ColumnWidth 12.5 is the COM value that this engine rounds down to the
#    narrower stored column width used in the updated report.
$newWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth   # column E
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth   # column F

$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth       # column C

$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth       # column C
